$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 18500
$ws.Range("J13").Value = 18500
$ws.Range("L13").Value = 18500
$ws.Range("N13").Value = -18838
$ws.Range("H40").Value = 38463620
$ws.Range("I40").Value = 1774.5
$ws.Range("J40").Value = 45456680
$ws.Range("K40").Value = 1774.5
$ws.Range("L40").Value = 45456680
$ws.Range("M40").Value = -1599.5
$ws.Range("N40").Value = -45457030
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H106").Value = 2782.652
$ws.Range("I106").Value = 2637.5625
$ws.Range("J106").Value = 3114.2856
$ws.Range("K106").Value = 2637.5625
$ws.Range("L106").Value = 3114.2856
$ws.Range("M106").Value = -2006.5625
$ws.Range("N106").Value = -4376.2856
$ws.Range("H137").Value = 971.4286
$ws.Range("I137").Value = 888.8889
$ws.Range("J137").Value = 1120
$ws.Range("K137").Value = 2666.6667
$ws.Range("L137").Value = 3360
$ws.Range("M137").Value = -116.6667000000002
$ws.Range("N137").Value = -8460
$ws.Range("H138").Value = 1528.2
$ws.Range("I138").Value = 693.03174
$ws.Range("J138").Value = 2950.2432
$ws.Range("K138").Value = 2079.09522
$ws.Range("L138").Value = 8850.729599999999
$ws.Range("M138").Value = 3060.90478
$ws.Range("N138").Value = -19130.7296
$ws.Range("H141").Value = 2825.0193
$ws.Range("I141").Value = 691.55554
$ws.Range("J141").Value = 7625.3125
$ws.Range("K141").Value = 2074.66662
$ws.Range("L141").Value = 22875.9375
$ws.Range("M141").Value = 3105.33338
$ws.Range("N141").Value = -33235.9375

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3765583.5
$ws.Range("I32").Value = 4649325
$ws.Range("J32").Value = 22677.117
$ws.Range("K32").Value = 4649325
$ws.Range("L32").Value = 22677.117
$ws.Range("M32").Value = -4649038
$ws.Range("N32").Value = -23251.117
$ws.Range("H74").Value = 804.58185
$ws.Range("I74").Value = 737.617
$ws.Range("K74").Value = 737.617
$ws.Range("M74").Value = 136.383
$ws.Range("H77").Value = 804.58185
$ws.Range("I77").Value = 737.617
$ws.Range("K77").Value = 3688.085
$ws.Range("M77").Value = 679.915
$ws.Range("H98").Value = 19051.666
$ws.Range("J98").Value = 19051.666
$ws.Range("L98").Value = 19051.666
$ws.Range("N98").Value = -25041.666
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1812.9302
$ws.Range("I86").Value = 1665.9395
$ws.Range("K86").Value = 1665.9395
$ws.Range("M86").Value = -542.9395
$ws.Range("H89").Value = 1812.9302
$ws.Range("I89").Value = 1665.9395
$ws.Range("K89").Value = 8329.6975
$ws.Range("M89").Value = -2713.6975
$ws.Range("H102").Value = 24777.75
$ws.Range("I102").Value = 10500
$ws.Range("J102").Value = 39055.5
$ws.Range("K102").Value = 10500
$ws.Range("L102").Value = 39055.5
$ws.Range("M102").Value = -7255
$ws.Range("N102").Value = -45545.5
$ws.Range("H105").Value = 5424.5
$ws.Range("I105").Value = 6186.6665
$ws.Range("J105").Value = 4800.909
$ws.Range("K105").Value = 6186.6665
$ws.Range("L105").Value = 4800.909
$ws.Range("M105").Value = -4439.6665
$ws.Range("N105").Value = -8294.909

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 600
$ws.Range("I17").Value = 600
$ws.Range("K17").Value = 600
$ws.Range("M17").Value = -426
$ws.Range("H31").Value = 2582.6858
$ws.Range("I31").Value = 2599.8235
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 2599.8235
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = -2304.8235
$ws.Range("N31").Value = -2590
$ws.Range("H34").Value = 2582.6858
$ws.Range("I34").Value = 2599.8235
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 2599.8235
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = -2397.8235
$ws.Range("N34").Value = -2404
$ws.Range("H134").Value = 1944.3334
$ws.Range("I134").Value = 1844.4
$ws.Range("K134").Value = 5533.200000000001
$ws.Range("M134").Value = -2998.200000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 538.46155
$ws.Range("I21").Value = 400
$ws.Range("J21").Value = 760
$ws.Range("K21").Value = 1200
$ws.Range("L21").Value = 2280
$ws.Range("M21").Value = -1027
$ws.Range("N21").Value = -2626
$ws.Range("H122").Value = 715058.5600000001
$ws.Range("J122").Value = 1000908
$ws.Range("L122").Value = 9008172
$ws.Range("N122").Value = -9013072
$ws.Range("H131").Value = 824.0599999999999
$ws.Range("I131").Value = 401.625
$ws.Range("J131").Value = 860.79346
$ws.Range("K131").Value = 1204.875
$ws.Range("L131").Value = 2582.38038
$ws.Range("M131").Value = 3835.125
$ws.Range("N131").Value = -12662.38038

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3824.718
$ws.Range("I70").Value = 3611.7742
$ws.Range("K70").Value = 3611.7742
$ws.Range("M70").Value = -3341.7742
$ws.Range("H73").Value = 3824.718
$ws.Range("I73").Value = 3611.7742
$ws.Range("K73").Value = 3611.7742
$ws.Range("M73").Value = -2675.7742
$ws.Range("H126").Value = 3432.6
$ws.Range("I126").Value = 3518.6667
$ws.Range("J126").Value = 3303.5
$ws.Range("K126").Value = 10556.0001
$ws.Range("L126").Value = 9910.5
$ws.Range("M126").Value = -8086.000100000001
$ws.Range("N126").Value = -14850.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 843534.5
$ws.Range("I40").Value = 1011841.4
$ws.Range("K40").Value = 1011841.4
$ws.Range("M40").Value = -1011705.4
$ws.Range("H46").Value = 1183.4
$ws.Range("I46").Value = 1154.25
$ws.Range("K46").Value = 1154.25
$ws.Range("M46").Value = -966.25
$ws.Range("H132").Value = 2179.178
$ws.Range("I132").Value = 1883.0212
$ws.Range("J132").Value = 2714.5386
$ws.Range("K132").Value = 5649.063599999999
$ws.Range("L132").Value = 8143.6158
$ws.Range("M132").Value = -3119.063599999999
$ws.Range("N132").Value = -13203.6158
$ws.Range("H133").Value = 42432.6
$ws.Range("J133").Value = 42432.6
$ws.Range("L133").Value = 42432.6
$ws.Range("N133").Value = -47492.6
$ws.Range("H136").Value = 4474.067
$ws.Range("I136").Value = 1731.1
$ws.Range("J136").Value = 9960
$ws.Range("K136").Value = 5193.299999999999
$ws.Range("L136").Value = 29880
$ws.Range("M136").Value = -2643.299999999999
$ws.Range("N136").Value = -34980

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 855.53845
$ws.Range("I107").Value = 1091.3334
$ws.Range("J107").Value = 325
$ws.Range("K107").Value = 3274.0002
$ws.Range("L107").Value = 975
$ws.Range("M107").Value = -1354.0002
$ws.Range("N107").Value = -4815
$ws.Range("H132").Value = 1147.6274
$ws.Range("I132").Value = 743.8182
$ws.Range("J132").Value = 1887.9445
$ws.Range("K132").Value = 2231.4546
$ws.Range("L132").Value = 5663.833500000001
$ws.Range("M132").Value = 298.5454
$ws.Range("N132").Value = -10723.8335
$ws.Range("H136").Value = 828
$ws.Range("I136").Value = 556.5
$ws.Range("K136").Value = 1669.5
$ws.Range("M136").Value = 880.5
